# "make attunits available for same time as defunits"
#
# Techtree: research "Bombenteppich" (K-column) and "Tragetasche" (L-column)
# are replaced by two new researches "Plasmatechnik" / "Superkompression"
# whose cost/time rows in Forschungen_Gebäude are cut down to match the
# def-unit ("Laser" / "leichte waffentechnik") research times, so the
# attack-unit tech unlocks in the same time as the matching defense-unit
# tech. A couple of stray Techtree cells (which unit needs which building)
# get swapped around too, and the active sheet/selection moves on to the
# Forschungen_Gebäude sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Einheiten")
$ws2 = $wb.Worksheets.Item("Techtree")
$ws3 = $wb.Worksheets.Item("Forschungen_Gebäude")

# --- Techtree sheet -------------------------------------------------
# Rename the research headers in row 1: K1/L1 get new tech names.
$ws2.Range("K1").Value = "Plasmatechnik"
$ws2.Range("L1").Value = "Superkompression"

# Building-requirement cells that moved between columns.
$ws2.Range("H8").Value = "-"
$ws2.Range("K8").Value = 1

$ws2.Range("D9").Value = 1
$ws2.Range("H9").Value = 1
$ws2.Range("K9").Value = "-"

# Forschungslabor-time column: 10 -> 1 for the three entries that used it.
$ws2.Range("O17").Value = 1
$ws2.Range("O18").Value = 1
$ws2.Range("O19").Value = 1

# Move the Techtree sheet's saved selection.
[void]$ws2.Range("K9").Select()

# --- Forschungen_Gebäude sheet --------------------------------------
# Laser research (row 5): costs /10, time 5h -> 1h (same as Mili/row 11-ish ratio)
$ws3.Range("B5").Value = 3000
$ws3.Range("C5").Value = 2000
$ws3.Range("D5").Value = 1000
$ws3.Range("E5").Value = 1000
$ws3.Range("F5").Value = 1/24

# leichte waffentechnik research (row 6)
$ws3.Range("B6").Value = 2000
$ws3.Range("C6").Value = 2000
$ws3.Range("D6").Value = 2000
$ws3.Range("E6").Value = 1000
$ws3.Range("F6").Value = 1/24

# row 7 becomes the new "Plasmatechnik" research (was "Bombenteppich")
$ws3.Range("A7").Value = "Plasmatechnik"
$ws3.Range("B7").Value = 2000
$ws3.Range("C7").Value = 3000
$ws3.Range("D7").Value = 1000
$ws3.Range("E7").Value = 1000
$ws3.Range("F7").Value = 1/24

# row 8 becomes the new "Superkompression" research (was "Tragetasche")
$ws3.Range("A8").Value = "Superkompression"

# Move selection to A9 and make this the active/visible sheet.
[void]$ws3.Range("A9").Select()
[void]$ws3.Activate()
